$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '30.346.56'
$ws.Range('E2').Value = '  -2.75%  '
$ws.Range('D3').Value = '1.941.39'
$ws.Range('E3').Value = '  -2.80%  '
$ws.Range('D4').Value = '0.9996'
$ws.Range('E4').Value = '  -0.06%  '
$ws.Range('D5').Value = '251.84'
$ws.Range('E5').Value = '  -1.65%  '
$ws.Range('D6').Value = '0.7202'
$ws.Range('E6').Value = '  -7.67%  '
$ws.Range('D7').Value = '0.9996'
$ws.Range('E7').Value = '  +0.00%  '
$ws.Range('D8').Value = '0.3349'
$ws.Range('E8').Value = '  -3.75%  '
$ws.Range('D9').Value = '28.80'
$ws.Range('E9').Value = '  +1.99%  '
$ws.Range('D10').Value = '0.07384'
$ws.Range('E10').Value = '  +5.50%  '
$ws.Range('D11').Value = '0.8184'
$ws.Range('E11').Value = '  -3.93%  '
$ws.Range('E12').Value = '  -0.42%  '
$ws.Range('D13').Value = '1.940.75'
$ws.Range('E13').Value = '  -2.81%  '
$ws.Range('D15').Value = '95.46'
$ws.Range('E15').Value = '  -5.07%  '
$ws.Range('D16').Value = '14.95'
$ws.Range('E16').Value = '  -3.19%  '
$ws.Range('D17').Value = '30.372.17'
$ws.Range('E17').Value = '  -2.72%  '
$ws.Range('D18').Value = '0.000008367'
$ws.Range('E18').Value = '  +5.77%  '
$ws.Range('D19').Value = '254.65'
$ws.Range('E19').Value = '  -7.10%  '
$ws.Range('D20').Value = '5.880'
$ws.Range('E20').Value = '  -0.14%  '
$ws.Range('D21').Value = '2.196.11'
$ws.Range('E21').Value = '  -2.80%  '
$ws.Range('D22').Value = '0.9998'
$ws.Range('E22').Value = '  -0.02%  '
$ws.Range('D23').Value = '0.9991'
$ws.Range('E23').Value = '  -0.11%  '
$ws.Range('D24').Value = '7.001'
$ws.Range('E24').Value = '  -1.47%  '
$ws.Range('D25').Value = '9.903'
$ws.Range('E25').Value = '  -1.47%  '
$ws.Range('D26').Value = '161.47'
$ws.Range('E26').Value = '  -1.81%  '
$ws.Range('D27').Value = '2.444'
$ws.Range('E27').Value = '  +5.30%  '
$ws.Range('D28').Value = '19.43'
$ws.Range('E28').Value = '  -2.40%  '
$ws.Range('D29').Value = '0.1316'
$ws.Range('E29').Value = '  -11.04%  '
$ws.Range('D30').Value = '1.573'
$ws.Range('E30').Value = '  -2.26%  '
$ws.Range('E31').Value = '  -0.93%  '
$ws.Range('D32').Value = '4.495'
$ws.Range('E32').Value = '  -2.45%  '
$ws.Range('D33').Value = '4.280'
$ws.Range('E33').Value = '  -2.73%  '
$ws.Range('D34').Value = '0.05317'
$ws.Range('E34').Value = '  +1.65%  '
$ws.Range('D35').Value = '1.323'
$ws.Range('E35').Value = '  +7.04%  '
$ws.Range('D36').Value = '0.7649'
$ws.Range('E36').Value = '  -1.45%  '
$ws.Range('D37').Value = '2.752'
$ws.Range('E37').Value = '  -0.16%  '
$ws.Range('D38').Value = '0.02002'
$ws.Range('E38').Value = '  +0.10%  '
$ws.Range('D39').Value = '2.846'
$ws.Range('E39').Value = '  -1.91%  '
$ws.Range('D40').Value = '81.50'
$ws.Range('E40').Value = '  +2.67%  '
$ws.Range('E41').Value = '  -1.06%  '
$ws.Range('D42').Value = '0.4584'
$ws.Range('E42').Value = '  -1.91%  '
$ws.Range('D43').Value = '2.044'
$ws.Range('E43').Value = '  -4.58%  '
$ws.Range('D44').Value = '0.8478'
$ws.Range('E44').Value = '  -0.32%  '
$ws.Range('D45').Value = '1.0000'
$ws.Range('E45').Value = '  +0.05%  '
$ws.Range('D46').Value = '103.30'
$ws.Range('E46').Value = '  -2.31%  '
$ws.Range('D47').Value = '9.915'
$ws.Range('E47').Value = '  +0.07%  '
$ws.Range('D48').Value = '7.509'
$ws.Range('E48').Value = '  -2.62%  '
$ws.Range('D49').Value = '37.33'
$ws.Range('E49').Value = '  +1.32%  '
$ws.Range('D50').Value = '0.4215'
$ws.Range('E50').Value = '  -2.50%  '

$ws.Range('B51').Value = 'NEARProtocol'
$ws.Range('C51').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D51').Value = '1.513'
$ws.Range('E51').Value = '  -1.95%  '
